$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update Maximo (C2) ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 696.198805276246

# --- Sheet "Solucion": update Pedido/Salida assignment (mutated solution) ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$solucionData = @(
    @("Pedido_19", "S001"),
    @("Pedido_50", "S025"),
    @("Pedido_23", "S041"),
    @("Pedido_2", "S065"),
    @("Pedido_11", "S029"),
    @("Pedido_24", "S005"),
    @("Pedido_48", "S045"),
    @("Pedido_32", "S069"),
    @("Pedido_37", "S026"),
    @("Pedido_31", "S002"),
    @("Pedido_22", "S042"),
    @("Pedido_38", "S066"),
    @("Pedido_73", "S046"),
    @("Pedido_36", "S070"),
    @("Pedido_6", "S006"),
    @("Pedido_12", "S030"),
    @("Pedido_47", "S043"),
    @("Pedido_59", "S003"),
    @("Pedido_13", "S067"),
    @("Pedido_33", "S027"),
    @("Pedido_10", "S047"),
    @("Pedido_3", "S071"),
    @("Pedido_60", "S007"),
    @("Pedido_35", "S031"),
    @("Pedido_7", "S044"),
    @("Pedido_9", "S004"),
    @("Pedido_34", "S028"),
    @("Pedido_76", "S048"),
    @("Pedido_74", "S068"),
    @("Pedido_25", "S008"),
    @("Pedido_53", "S049"),
    @("Pedido_40", "S009"),
    @("Pedido_71", "S032"),
    @("Pedido_75", "S072"),
    @("Pedido_65", "S053"),
    @("Pedido_67", "S013"),
    @("Pedido_51", "S073"),
    @("Pedido_78", "S033"),
    @("Pedido_46", "S050"),
    @("Pedido_29", "S010"),
    @("Pedido_43", "S077"),
    @("Pedido_44", "S037"),
    @("Pedido_55", "S054"),
    @("Pedido_17", "S014"),
    @("Pedido_64", "S034"),
    @("Pedido_15", "S074"),
    @("Pedido_4", "S051"),
    @("Pedido_1", "S011"),
    @("Pedido_58", "S038"),
    @("Pedido_27", "S078"),
    @("Pedido_30", "S015"),
    @("Pedido_41", "S055"),
    @("Pedido_45", "S035"),
    @("Pedido_61", "S075"),
    @("Pedido_16", "S012"),
    @("Pedido_28", "S052"),
    @("Pedido_26", "S079"),
    @("Pedido_5", "S016"),
    @("Pedido_49", "S039"),
    @("Pedido_8", "S056"),
    @("Pedido_21", "S076"),
    @("Pedido_69", "S036"),
    @("Pedido_72", "S017"),
    @("Pedido_18", "S040"),
    @("Pedido_62", "S021"),
    @("Pedido_14", "S057"),
    @("Pedido_52", "S080"),
    @("Pedido_70", "S018"),
    @("Pedido_63", "S061"),
    @("Pedido_77", "S022"),
    @("Pedido_80", "S058"),
    @("Pedido_54", "S019"),
    @("Pedido_68", "S062"),
    @("Pedido_20", "S059"),
    @("Pedido_66", "S023"),
    @("Pedido_56", "S020"),
    @("Pedido_39", "S063"),
    @("Pedido_79", "S024"),
    @("Pedido_57", "S060"),
    @("Pedido_42", "S064")
)
for ($i = 0; $i -lt $solucionData.Count; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $solucionData[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $solucionData[$i][1]
}

# --- Sheet "Metricas": update Tiempo (Wmax) per zone ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 677.2111579630231
$wsMetricas.Range("B3").Value = 467.3903124662126
$wsMetricas.Range("B4").Value = 696.198805276246
$wsMetricas.Range("B5").Value = 487.8754189642124
